$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4272.923
$ws.Range("J64").Value = 4664.143
$ws.Range("L64").Value = 4664.143
$ws.Range("N64").Value = -5160.143
$ws.Range("H67").Value = 4272.923
$ws.Range("J67").Value = 4664.143
$ws.Range("L67").Value = 4664.143
$ws.Range("N67").Value = -6380.143
$ws.Range("H74").Value = 2224.2415
$ws.Range("I74").Value = 1659.2273
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 1659.2273
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -723.2273
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 2224.2415
$ws.Range("I77").Value = 1659.2273
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 8296.136500000001
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -3616.136500000001
$ws.Range("N77").Value = -29360
$ws.Range("H106").Value = 1494.1666
$ws.Range("I106").Value = 1311.3636
$ws.Range("K106").Value = 1311.3636
$ws.Range("M106").Value = -680.3635999999999
$ws.Range("H138").Value = 3015.8333
$ws.Range("J138").Value = 3488.1282
$ws.Range("L138").Value = 10464.3846
$ws.Range("N138").Value = -20744.3846
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7303.0835
$ws.Range("I61").Value = 14241.333
$ws.Range("J61").Value = 4990.3335
$ws.Range("K61").Value = 14241.333
$ws.Range("L61").Value = 4990.3335
$ws.Range("M61").Value = -14029.333
$ws.Range("N61").Value = -5414.3335
$ws.Range("H74").Value = 890.7857
$ws.Range("I74").Value = 544.3333
$ws.Range("J74").Value = 1150.625
$ws.Range("K74").Value = 544.3333
$ws.Range("L74").Value = 1150.625
$ws.Range("M74").Value = 329.6667
$ws.Range("N74").Value = -2898.625
$ws.Range("H77").Value = 890.7857
$ws.Range("I77").Value = 544.3333
$ws.Range("J77").Value = 1150.625
$ws.Range("K77").Value = 2721.6665
$ws.Range("L77").Value = 5753.125
$ws.Range("M77").Value = 1646.3335
$ws.Range("N77").Value = -14489.125
$ws.Range("H88").Value = 50445.43
$ws.Range("I88").Value = 1581.2
$ws.Range("J88").Value = 65715.5
$ws.Range("K88").Value = 1581.2
$ws.Range("L88").Value = 65715.5
$ws.Range("M88").Value = -1175.2
$ws.Range("N88").Value = -66527.5
$ws.Range("H91").Value = 50445.43
$ws.Range("I91").Value = 1581.2
$ws.Range("J91").Value = 65715.5
$ws.Range("K91").Value = 1581.2
$ws.Range("L91").Value = 65715.5
$ws.Range("M91").Value = -177.2
$ws.Range("N91").Value = -68523.5
$ws.Range("H110").Value = 2003.5294
$ws.Range("I110").Value = 1788.8889
$ws.Range("J110").Value = 2245
$ws.Range("K110").Value = 1788.8889
$ws.Range("L110").Value = 2245
$ws.Range("M110").Value = 256.1111000000001
$ws.Range("N110").Value = -6335
$ws.Range("H132").Value = 59023.555
$ws.Range("I132").Value = 10012
$ws.Range("J132").Value = 65150
$ws.Range("K132").Value = 30036
$ws.Range("L132").Value = 195450
$ws.Range("M132").Value = -27506
$ws.Range("N132").Value = -200510
$ws.Range("H136").Value = 7303.0835
$ws.Range("I136").Value = 14241.333
$ws.Range("J136").Value = 4990.3335
$ws.Range("K136").Value = 42723.999
$ws.Range("L136").Value = 14971.0005
$ws.Range("M136").Value = -40173.999
$ws.Range("N136").Value = -20071.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 312.09525
$ws.Range("I22").Value = 320.6
$ws.Range("J22").Value = 142
$ws.Range("K22").Value = 320.6
$ws.Range("L22").Value = 142
$ws.Range("M22").Value = -147.6
$ws.Range("N22").Value = -488
$ws.Range("H134").Value = 34863.258
$ws.Range("I134").Value = 34863.258
$ws.Range("K134").Value = 104589.774
$ws.Range("M134").Value = -102054.774
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18154.477
$ws.Range("I31").Value = 66300.2
$ws.Range("J31").Value = 3108.9375
$ws.Range("K31").Value = 66300.2
$ws.Range("L31").Value = 3108.9375
$ws.Range("M31").Value = -66005.2
$ws.Range("N31").Value = -3698.9375
$ws.Range("H34").Value = 18154.477
$ws.Range("I34").Value = 66300.2
$ws.Range("J34").Value = 3108.9375
$ws.Range("K34").Value = 66300.2
$ws.Range("L34").Value = 3108.9375
$ws.Range("M34").Value = -66098.2
$ws.Range("N34").Value = -3512.9375
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3860.258
$ws.Range("J107").Value = 741.2222
$ws.Range("L107").Value = 2223.6666
$ws.Range("N107").Value = -6063.6666
$ws.Range("H131").Value = 176287.95
$ws.Range("J131").Value = 193144.6
$ws.Range("L131").Value = 579433.8
$ws.Range("N131").Value = -589513.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8285.789000000001
$ws.Range("I80").Value = 11240.909
$ws.Range("J80").Value = 4222.5
$ws.Range("K80").Value = 11240.909
$ws.Range("L80").Value = 4222.5
$ws.Range("M80").Value = -10242.909
$ws.Range("N80").Value = -6218.5
$ws.Range("H83").Value = 8285.789000000001
$ws.Range("I83").Value = 11240.909
$ws.Range("J83").Value = 4222.5
$ws.Range("K83").Value = 56204.545
$ws.Range("L83").Value = 21112.5
$ws.Range("M83").Value = -51212.545
$ws.Range("N83").Value = -31096.5
$ws.Range("H122").Value = 3649
$ws.Range("I122").Value = 2831.4443
$ws.Range("K122").Value = 8494.332900000001
$ws.Range("M122").Value = -6044.332900000001
$ws.Range("H132").Value = 221428.72
$ws.Range("I132").Value = 208800.2
$ws.Range("J132").Value = 253000
$ws.Range("K132").Value = 626400.6000000001
$ws.Range("L132").Value = 759000
$ws.Range("M132").Value = -623870.6000000001
$ws.Range("N132").Value = -764060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2407.389
$ws.Range("I93").Value = 2175.2666
$ws.Range("J93").Value = 3568
$ws.Range("K93").Value = 2175.2666
$ws.Range("L93").Value = 3568
$ws.Range("M93").Value = -927.2665999999999
$ws.Range("N93").Value = -6064
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1219.1666
$ws.Range("I81").Value = 1422.6
$ws.Range("J81").Value = 202
$ws.Range("K81").Value = 2845.2
$ws.Range("L81").Value = 404
$ws.Range("M81").Value = -1784.2
$ws.Range("N81").Value = -2526
$ws.Range("H84").Value = 1219.1666
$ws.Range("I84").Value = 1422.6
$ws.Range("J84").Value = 202
$ws.Range("K84").Value = 14226
$ws.Range("L84").Value = 2020
$ws.Range("M84").Value = -8922
$ws.Range("N84").Value = -12628
$ws.Range("H96").Value = 3800.0667
$ws.Range("I96").Value = 1916.8334
$ws.Range("J96").Value = 5055.5557
$ws.Range("K96").Value = 1916.8334
$ws.Range("L96").Value = 5055.5557
$ws.Range("M96").Value = -543.8334
$ws.Range("N96").Value = -7801.5557

Write-Host "Applied all changes"